$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (number format, font, border, alignment) from the last
# existing data row (A269) so the new date cells in column A match the
# existing formatting (style index reused, no new style entries created).
$ws.Range("A269").Copy() | Out-Null

$startRow = 270
$endRow = 301
$startSerial = 44344

for ($r = $startRow; $r -le $endRow; $r++) {
    $serial = $startSerial + ($r - $startRow)

    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = $serial
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0
}

$excel.CutCopyMode = 0
